$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original formatting of the Price column while writing values as text,
# so numeric-looking strings (e.g. "1.002") are not auto-converted to numbers.
$origPriceStyle = $ws.Range("D2:D51").Style
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "24.668.00"
$ws.Range("E2").Value = "  +0.15%  "

# Row 3
$ws.Range("D3").Value = "1.690.04"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "315.50"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("E6").Value = "  +0.23%  "

# Row 7
$ws.Range("D7").Value = "0.3945"
$ws.Range("E7").Value = "  -0.61%  "

# Row 8
$ws.Range("D8").Value = "0.4055"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "1.487"
$ws.Range("E9").Value = "  -1.86%  "

# Row 10
$ws.Range("D10").Value = "1.003"
$ws.Range("E10").Value = "  +0.17%  "

# Row 11
$ws.Range("D11").Value = "52.46"
$ws.Range("E11").Value = "  -2.31%  "

# Row 12
$ws.Range("D12").Value = "0.08841"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("D13").Value = "7.243"
$ws.Range("E13").Value = "  -1.26%  "

# Row 14
$ws.Range("D14").Value = "23.54"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15
$ws.Range("D15").Value = "8.044"
$ws.Range("E15").Value = "  +6.83%  "

# Row 16
$ws.Range("D16").Value = "0.00001316"
$ws.Range("E16").Value = "  -0.84%  "

# Row 17
$ws.Range("D17").Value = "1.693.05"
$ws.Range("E17").Value = "  -0.40%  "

# Row 18
$ws.Range("D18").Value = "99.58"
$ws.Range("E18").Value = "  -1.43%  "

# Row 19
$ws.Range("D19").Value = "0.07021"
$ws.Range("E19").Value = "  -1.31%  "

# Row 20
$ws.Range("D20").Value = "19.51"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("D21").Value = "6.991"
$ws.Range("E21").Value = "  +3.29%  "

# Row 22
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("D23").Value = "14.31"
$ws.Range("E23").Value = "  +1.31%  "

# Row 24
$ws.Range("D24").Value = "24.664.94"
$ws.Range("E24").Value = "  +0.17%  "

# Row 25
$ws.Range("D25").Value = "3.308"
$ws.Range("E25").Value = "  +9.71%  "

# Row 26
$ws.Range("E26").Value = "  +1.88%  "

# Row 27
$ws.Range("E27").Value = "  +1.14%  "

# Row 28
$ws.Range("D28").Value = "162.55"
$ws.Range("E28").Value = "  +1.98%  "

# Row 29
$ws.Range("D29").Value = "135.53"
$ws.Range("E29").Value = "  +1.08%  "

# Row 30
$ws.Range("D30").Value = "5.175"

# Row 31
$ws.Range("D31").Value = "7.593"
$ws.Range("E31").Value = "  +2.93%  "

# Row 32
$ws.Range("D32").Value = "1.878.80"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.08545"
$ws.Range("E33").Value = "  -1.80%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.057"
$ws.Range("E34").Value = "  -3.23%  "

# Row 35
$ws.Range("D35").Value = "7.086"
$ws.Range("E35").Value = "  -3.21%  "

# Row 36
$ws.Range("D36").Value = "11.21"
$ws.Range("E36").Value = "  +0.96%  "

# Row 37
$ws.Range("D37").Value = "0.2734"
$ws.Range("E37").Value = "  +0.27%  "

# Row 38
$ws.Range("D38").Value = "1.882"
$ws.Range("E38").Value = "  -4.17%  "

# Row 39
$ws.Range("D39").Value = "14.47"
$ws.Range("E39").Value = "  -2.06%  "

# Row 40
$ws.Range("D40").Value = "0.09198"
$ws.Range("E40").Value = "  +2.09%  "

# Row 41
$ws.Range("D41").Value = "0.02720"
$ws.Range("E41").Value = "  -1.97%  "

# Row 42
$ws.Range("D42").Value = "1.463"
$ws.Range("E42").Value = "  -1.26%  "

# Row 43
$ws.Range("D43").Value = "0.7630"
$ws.Range("E43").Value = "  -0.79%  "

# Row 44
$ws.Range("D44").Value = "16.13"
$ws.Range("E44").Value = "  +2.89%  "

# Row 45
$ws.Range("D45").Value = "2.587"
$ws.Range("E45").Value = "  +5.17%  "

# Row 46
$ws.Range("D46").Value = "0.7130"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47
$ws.Range("D47").Value = "4.204"
$ws.Range("E47").Value = "  +0.57%  "

# Row 48
$ws.Range("E48").Value = "  +0.26%  "

# Row 49
$ws.Range("D49").Value = "140.12"
$ws.Range("E49").Value = "  -0.76%  "

# Row 50
$ws.Range("E50").Value = "  +1.43%  "

# Row 51
$ws.Range("E51").Value = "  -0.32%  "

# Restore the original formatting on the Price column
$ws.Range("D2:D51").Style = $origPriceStyle
